# Add 12 new daily columns (LR:MC) to the "mobility" sheet, continuing the
# date series in row 1 and the data series in rows 2-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dates = @(44166,44167,44168,44169,44170,44171,44172,44173,44174,44175,44176,44177)
$row2  = @(60, 57.87, 57.62, 61.4, 57.16, 47.64, 60.66, 58.97, 59.2, 59.07, 60.04, 59.21)
$row3  = @(44.94, 44.16, 45.75, 48.53, 48.42, 44.62, 44.48, 43.82, 43.81, 46.79, 53.32, 53.28)
$row4  = @(62.65, 61.41, 62.8, 60.66, 58.4, 49.42, 60.38, 60.58, 60.55, 54.09, 63.32, 59.22)
$row5  = @(31.9, 32.47, 31.45, 31.9, 28.83, 23.95, 32.630000000000003, 30.04, 30.49, 29.06, 29.79, 25.19)

# Row 1, columns G..LQ (7..329) already hold dates formatted as dates (style "1").
# New columns start right after LQ, i.e. column 330 (LR) through column 341 (MC).
$startCol = 330

for ($i = 0; $i -lt $dates.Length; $i++) {
    $col = $startCol + $i

    $cell1 = $ws.Cells.Item(1, $col)
    $cell1.Value = $dates[$i]
    $cell1.NumberFormat = "m/d/yy"

    $ws.Cells.Item(2, $col).Value = $row2[$i]
    $ws.Cells.Item(3, $col).Value = $row3[$i]
    $ws.Cells.Item(4, $col).Value = $row4[$i]
    $ws.Cells.Item(5, $col).Value = $row5[$i]
}
